{"js": "// Strip the filled-in sample answers from the form's field labels, leaving\n// just the bold label text (e.g. \"Jm\u00e9no a p\u0159\u00edjmen\u00ed: Martin Vader: Martin Vader\"\n// -> \"Jm\u00e9no a p\u0159\u00edjmen\u00ed:\"). Some labels were already bold; others need <w:b/>\n// + <w:bCs/> added. The two section headings (C, D) just lose their\n// appended answer text, no bold change. The two \"Elektronick\u00fd podpis: \"\n// labels only lose their trailing space.\n\nasync function replaceLabel(body, oldText, newText, makeBold) {\n  const results = body.search(oldText, { matchCase: true });\n  await context.sync();\n\n  results.load(\"items\");\n  await context.sync();\n\n  results.items.forEach((range) => {\n    range.insertText(newText, Word.InsertLocation.replace);\n    if (makeBold) {\n      range.font.bold = true;\n      range.font.boldBidirectional = true;\n    }\n  });\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// Section A \u2014 \u00dadaje o zam\u011bstnanci\nawait replaceLabel(body, \"Jm\u00e9no a p\u0159\u00edjmen\u00ed: Martin Vader: Martin Vader\", \"Jm\u00e9no a p\u0159\u00edjmen\u00ed:\", true);\nawait replaceLabel(body, \"Osobn\u00ed \u010d\u00edslo: 2015: 2015\", \"Osobn\u00ed \u010d\u00edslo:\", true);\nawait replaceLabel(body, \"\u00datvar / odd\u011blen\u00ed: magiologie: \u00fatvar pro magiologii\", \"\u00datvar / odd\u011blen\u00ed:\", true);\nawait replaceLabel(body, \"Telefon / e\u2011mail: mkk@magik.cz\", \"Telefon / e\u2011mail:\", false);\n\n// Section B \u2014 \u00dadaje o cest\u011b\nawait replaceLabel(body, \"Datum a \u010das odjezdu: 30.11. (\u010das bude dopln\u011bn): 30.11.2025\", \"Datum a \u010das odjezdu:\", true);\nawait replaceLabel(body, \"Datum a \u010das n\u00e1vratu: 5.12. (\u010das bude dopln\u011bn): 5.12.2025\", \"Datum a \u010das n\u00e1vratu:\", true);\nawait replaceLabel(body, \"M\u00edsto kon\u00e1n\u00ed cesty (m\u011bsto, adresa): Fakultn\u00ed nemocnice u sv. Anny v Brn\u011b: Brno\", \"M\u00edsto kon\u00e1n\u00ed cesty (m\u011bsto, adresa):\", true);\nawait replaceLabel(body, \"\u00da\u010del cesty (stru\u010dn\u00fd popis): \u0161kolen\u00ed student\u016f o magii\", \"\u00da\u010del cesty (stru\u010dn\u00fd popis):\", false);\n\n// Section C \u2014 Zp\u016fsob dopravy (heading, no bold)\nawait replaceLabel(body, \"C \u2014 Zp\u016fsob dopravy (vyberte / ozna\u010dte): hromadn\u00e1 doprava \u2013 vlak\", \"C \u2014 Zp\u016fsob dopravy (vyberte / ozna\u010dte)\", false);\n\n// Section D \u2014 Odhadovan\u00e9 n\u00e1klady (heading, no bold)\nawait replaceLabel(body, \"D \u2014 Odhadovan\u00e9 n\u00e1klady (nepovinn\u00e9): Odhad n\u00e1klad\u016f: 80 K\u010d \u2013 j\u00edzdn\u00e9 vlakem.\", \"D \u2014 Odhadovan\u00e9 n\u00e1klady (nepovinn\u00e9)\", false);\n\n// Section E \u2014 Prohl\u00e1\u0161en\u00ed \u017eadatele\nawait replaceLabel(body, \"Datum: 27.11.2025\", \"Datum:\", false);\n\n// Both \"Elektronick\u00fd podpis: \" occurrences (sections E and F) \u2014 strip trailing space\nawait replaceLabel(body, \"Elektronick\u00fd podpis: \", \"Elektronick\u00fd podpis:\", false);\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Set-FieldLabel($OldText, $NewText, $MakeBold) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $found = $range.Find.Execute($OldText)\n    if ($found) {\n        $range.Text = $NewText\n        if ($MakeBold) {\n            $range.Font.Bold = 1\n            $range.Font.BoldBi = 1\n        }\n    }\n}\n\n# Section A \u2014 \u00dadaje o zam\u011bstnanci\nSet-FieldLabel \"Jm\u00e9no a p\u0159\u00edjmen\u00ed: Martin Vader: Martin Vader\" \"Jm\u00e9no a p\u0159\u00edjmen\u00ed:\" $true\nSet-FieldLabel \"Osobn\u00ed \u010d\u00edslo: 2015: 2015\" \"Osobn\u00ed \u010d\u00edslo:\" $true\nSet-FieldLabel \"\u00datvar / odd\u011blen\u00ed: magiologie: \u00fatvar pro magiologii\" \"\u00datvar / odd\u011blen\u00ed:\" $true\nSet-FieldLabel \"Telefon / e\u2011mail: mkk@magik.cz\" \"Telefon / e\u2011mail:\" $false\n\n# Section B \u2014 \u00dadaje o cest\u011b\nSet-FieldLabel \"Datum a \u010das odjezdu: 30.11. (\u010das bude dopln\u011bn): 30.11.2025\" \"Datum a \u010das odjezdu:\" $true\nSet-FieldLabel \"Datum a \u010das n\u00e1vratu: 5.12. (\u010das bude dopln\u011bn): 5.12.2025\" \"Datum a \u010das n\u00e1vratu:\" $true\nSet-FieldLabel \"M\u00edsto kon\u00e1n\u00ed cesty (m\u011bsto, adresa): Fakultn\u00ed nemocnice u sv. Anny v Brn\u011b: Brno\" \"M\u00edsto kon\u00e1n\u00ed cesty (m\u011bsto, adresa):\" $true\nSet-FieldLabel \"\u00da\u010del cesty (stru\u010dn\u00fd popis): \u0161kolen\u00ed student\u016f o magii\" \"\u00da\u010del cesty (stru\u010dn\u00fd popis):\" $false\n\n# Section C \u2014 Zp\u016fsob dopravy (heading)\nSet-FieldLabel \"C \u2014 Zp\u016fsob dopravy (vyberte / ozna\u010dte): hromadn\u00e1 doprava \u2013 vlak\" \"C \u2014 Zp\u016fsob dopravy (vyberte / ozna\u010dte)\" $false\n\n# Section D \u2014 Odhadovan\u00e9 n\u00e1klady (heading)\nSet-FieldLabel \"D \u2014 Odhadovan\u00e9 n\u00e1klady (nepovinn\u00e9): Odhad n\u00e1klad\u016f: 80 K\u010d \u2013 j\u00edzdn\u00e9 vlakem.\" \"D \u2014 Odhadovan\u00e9 n\u00e1klady (nepovinn\u00e9)\" $false\n\n# Section E \u2014 Prohl\u00e1\u0161en\u00ed \u017eadatele\nSet-FieldLabel \"Datum: 27.11.2025\" \"Datum:\" $false\n\n# Both \"Elektronick\u00fd podpis: \" occurrences (E and F sections) \u2014 strip trailing space\nfor ($i = 0; $i -lt 2; $i++) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Text = \"Elektronick\u00fd podpis: \"\n    $range.Find.MatchWildcards = $false\n    $found = $range.Find.Execute()\n    if ($found) {\n        $range.Text = \"Elektronick\u00fd podpis:\"\n    }\n}\n"}
